# Revert "Seasonality Index" (column L) values on the "Forecast Comparison"
# sheet back to their original figures (Optuna attempt rollback).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value  = 0.99
$ws.Range("L3").Value  = 1.08
$ws.Range("L4").Value  = 1.1
$ws.Range("L5").Value  = 1.19
$ws.Range("L6").Value  = 1
$ws.Range("L7").Value  = 1.04
$ws.Range("L8").Value  = 0.9399999999999999
$ws.Range("L10").Value = 0.93
$ws.Range("L11").Value = 0.96
$ws.Range("L12").Value = 0.99
$ws.Range("L13").Value = 0.97
$ws.Range("L14").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("L16").Value = 1.17
$ws.Range("L17").Value = 1.02
